$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 3 (a data row followed by its blank spacer row),
# pushing the existing log entries (and their spacer rows) down by two rows.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Fill in the new trade-plan entry for Friday 2017-01-06 in row 3.
$ws.Cells.Item(3, 1).Value = 20170105
$ws.Cells.Item(3, 2).Value = "Thursday"
$ws.Cells.Item(3, 3).Value = 20170106
$ws.Cells.Item(3, 4).Value = "Friday"
$ws.Cells.Item(3, 5).Value = "The report came out as -49, way less withdraw than the expected -72~-85. This is very bearish but NG quickly when it touched 3.170 and consolidated around 3.20, where we can feel the support. NG ended up a green candle with very bearish report, which simply means the bearish report were well expected and priced in already. But due to warm weather forecast, the NG continues to see presistence around 3.35, what's interesting is that NG continues to rebound after the outcry 14:30 but quickly drew back after. It's also interesting to notice the volume for DGAZ is relatively small while the volume of UGAZ is larger than normal, I have a strong feeling this could be a profit taking for NG shot and day trade opportunity to buy UGAZ tomorrow (or even hold over the weekend)"

$ws.Range("A3").WrapText = $true
$ws.Range("B3:G3").WrapText = $true
$ws.Range("B3:G3").HorizontalAlignment = -4131

$ws.Rows.Item(3).RowHeight = 115.2

$ws.Range("E3").Select()
